$wb = $excel.ActiveWorkbook

# --- Clean up the "房费" sheet: clear a handful of redundant empty/blank
# cells that carried a no-op border-less style left over from earlier
# edits (matches upstream author's save) ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B3:C3").ClearFormats()
$ws1.Range("D3:E3").Clear()

$ws1.Range("B7:E7").ClearFormats()
$ws1.Range("F7").Clear()

$ws1.Range("D11:E11").ClearFormats()
$ws1.Range("B11:C11").Clear()
$ws1.Range("F11").Clear()

$ws1.Range("E15:G15").ClearFormats()
$ws1.Range("B15:D15").Clear()
$ws1.Range("H15").Clear()

$ws1.Range("G19").Select() | Out-Null

# --- Add a new "Yang" worksheet at the end of the workbook with a small
# account-info table, and make it the active tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Yang"

$newSheet.Range("C2").Value = "aosiwei"
$newSheet.Range("D2").Value = "130plz"
$newSheet.Range("C3").Value = "mine"
$newSheet.Range("D3").Value = 121

$newSheet.Activate()
$newSheet.Range("D2").Select() | Out-Null
